$wb = $excel.ActiveWorkbook

# --- 1. Text update: "Ready for handoff" -> "In Translation" ---
# This string is used as the localization "Status" value on every sheet in
# the workbook (Overview's per-language status columns, and the Status
# column on each per-language sheet).
foreach ($ws in $wb.Worksheets) {
    $used = $ws.UsedRange
    foreach ($cell in $used.Cells) {
        # NB: keep the string literal on the left of -eq; some cells hold a
        # real Boolean (e.g. "True"/"False" columns) and PowerShell's -eq
        # coerces the right operand to the left operand's type, so
        # "$cell.Value() -eq 'Ready for handoff'" would wrongly match any
        # truthy Boolean cell.
        if ("Ready for handoff" -eq $cell.Value()) {
            $cell.Value = "In Translation"
        }
    }
}

# --- 2. Column width updates ---
# Overview sheet: columns E and F narrowed
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("E1").ColumnWidth = 12.576851254417766
$wsOverview.Range("F1").ColumnWidth = 12.576851254417766

# zh-cn sheet: column C narrowed
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("C1").ColumnWidth = 12.576851254417766

# de-de sheet: column C narrowed
$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("C1").ColumnWidth = 12.576851254417766
